$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix classification report sort order: the rating rows were written in the
# wrong order (A before AAA, B before BBB, C before CCC). Swap the data for
# row pairs (2,4), (5,7) and (8,10) so ratings run best-to-worst
# (AAA, AA, A, BBB, BB, B, CCC, CC, C, D). Rows 3, 6, 9 and 11 already hold
# the correct rating and are left untouched.
# Number-like values (precision/recall/F1/support) are stored as text in
# this report, so force text formatting before writing them back to avoid
# Excel re-interpreting them as numbers.

$changedRows = @(2, 4, 5, 7, 8, 10)
foreach ($r in $changedRows) {
    $ws.Range("B$r`:E$r").NumberFormat = "@"
}

$ws.Range("A2").Value = "AAA"
$ws.Range("B2").Value = "0.8519"
$ws.Range("C2").Value = "0.9583"
$ws.Range("D2").Value = "0.9020"
$ws.Range("E2").Value = "24"

$ws.Range("A4").Value = "A"
$ws.Range("B4").Value = "0.9604"
$ws.Range("C4").Value = "0.9327"
$ws.Range("D4").Value = "0.9463"
$ws.Range("E4").Value = "208"

$ws.Range("A5").Value = "BBB"
$ws.Range("B5").Value = "0.9669"
$ws.Range("C5").Value = "0.9669"
$ws.Range("D5").Value = "0.9669"
$ws.Range("E5").Value = "363"

$ws.Range("A7").Value = "B"
$ws.Range("B7").Value = "0.9359"
$ws.Range("C7").Value = "0.9481"
$ws.Range("D7").Value = "0.9419"
$ws.Range("E7").Value = "154"

$ws.Range("A8").Value = "CCC"
$ws.Range("B8").Value = "0.8065"
$ws.Range("C8").Value = "0.9615"
$ws.Range("D8").Value = "0.8772"
$ws.Range("E8").Value = "26"

$ws.Range("A10").Value = "C"
$ws.Range("B10").Value = "1.0000"
$ws.Range("C10").Value = "1.0000"
$ws.Range("D10").Value = "1.0000"
$ws.Range("E10").Value = "4"
